$d = $word.ActiveDocument

$old = "Datas das campanhas de 2018 que usam Perseu: 30 de outubro a 8 de novembro e 29 de novembro a 8 de dezembro"
$new = "Datas das campanhas de Hercules: 13-22 de junho, 12-21 de julho, 10-19 de agosto"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = $old
$find.Forward = $true
$find.Wrap = 1

$count = 0
while ($find.Execute()) {
    $rng = $find.Parent
    $rng.Delete()
    $rng.InsertAfter($new)
    $count = $count + 1
    if ($count -gt 20) { break }
}
Write-Output "Replaced $count occurrence(s)"
